$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 36.81180933333333
$ws.Range("N2").Value = 110.435428
$ws.Range("O2").Value = 0.2598784967371026
$ws.Range("P2").Value = 0.2598784967371026
$ws.Range("Q2").Value = 9.717875922287998
$ws.Range("R2").Value = 87.46088330059199
$ws.Range("S2").Value = 0.2598784967371026
$ws.Range("T2").Value = 0.2598784967371026

# Row 3
$ws.Range("O3").Value = 0.1970278712683331
$ws.Range("P3").Value = 0.197027871268333
$ws.Range("S3").Value = 0.1970278712683331
$ws.Range("T3").Value = 0.197027871268333

# Row 4
$ws.Range("M4").Value = 21.95609833333333
$ws.Range("N4").Value = 65.868295
$ws.Range("O4").Value = 0.1550023737603119
$ws.Range("P4").Value = 0.1550023737603119
$ws.Range("Q4").Value = 5.79614648682
$ws.Range("R4").Value = 52.16531838138
$ws.Range("S4").Value = 0.1550023737603119
$ws.Range("T4").Value = 0.1550023737603119

# Row 5
$ws.Range("M5").Value = 13.23098133333333
$ws.Range("N5").Value = 39.692944
$ws.Range("O5").Value = 0.09340609987756826
$ws.Range("P5").Value = 0.09340609987756825
$ws.Range("Q5").Value = 3.492820300224
$ws.Range("R5").Value = 31.435382702016
$ws.Range("S5").Value = 0.09340609987756826
$ws.Range("T5").Value = 0.09340609987756825

# Row 6
$ws.Range("M6").Value = 22.080681
$ws.Range("N6").Value = 66.242043
$ws.Range("O6").Value = 0.1558818838066577
$ws.Range("P6").Value = 0.1558818838066577
$ws.Range("Q6").Value = 5.829034815828
$ws.Range("R6").Value = 52.46131334245199
$ws.Range("S6").Value = 0.1558818838066577
$ws.Range("T6").Value = 0.1558818838066577

# Row 7
$ws.Range("M7").Value = 19.66149466666667
$ws.Range("N7").Value = 58.984484
$ws.Range("O7").Value = 0.1388032745500265
$ws.Range("P7").Value = 0.1388032745500265
$ws.Range("Q7").Value = 5.190398654063999
$ws.Range("R7").Value = 46.713587886576
$ws.Range("S7").Value = 0.1388032745500265
$ws.Range("T7").Value = 0.1388032745500265
